# [CC] correction of the reference values
#
# The roughness reference table (fmod in column A, R in column B,
# starting at row 3) is recomputed with corrected values and grows from
# 21 rows (3:23) to 24 rows (3:26). Column A is written in full first,
# then column B, so the corrected fmod values occupy the shared-string
# table before the corrected R values (matching how the sheet was
# originally authored).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmod = @(
    "11.087001", "13.655874", "17.425774", "22.323097", "28.37367",
    "34.94905",  "42.370758", "49.96757",  "56.425716", "63.209957",
    "69.42832",  "74.194626", "81.15967",  "87.05505",  "95.19742",
    "103.28752", "112.92756", "137.10373", "160.75127", "193.66571",
    "220.964",   "246.33324", "284.3601",  "332.1115"
)

$r = @(
    "0.10113691", "0.13728729", "0.19609652", "0.27899456", "0.39539358",
    "0.5409575",  "0.7005437",  "0.8453598",  "0.9321353",  "0.98051274",
    "0.99174285", "0.9837013",  "0.9491586",  "0.90162235", "0.80737907",
    "0.7173512",  "0.6151961",  "0.44188514", "0.34071133", "0.25057378",
    "0.19783457", "0.16634719", "0.13030091", "0.10086449"
)

$firstRow = 3

# Write column A (fmod) top to bottom ...
for ($i = 0; $i -lt $fmod.Length; $i++) {
    $row = $firstRow + $i
    $cell = $ws.Range("A$row")
    # Leading apostrophe forces text storage (matching the workbook's
    # existing text-typed reference values) without altering cell style.
    $cell.Value = "'" + $fmod[$i]
    $cell.Style = "Normal"
}

# ... then column B (R), so new unique strings are appended in the same
# order as in the source workbook (all fmod values, then all R values).
for ($i = 0; $i -lt $r.Length; $i++) {
    $row = $firstRow + $i
    $cell = $ws.Range("B$row")
    $cell.Value = "'" + $r[$i]
    $cell.Style = "Normal"
}

$lastRow = $firstRow + $fmod.Length - 1
[void]$ws.Range("B$firstRow`:B$lastRow").Select()
